# Auto-applies the cryptos.xlsx diff: updates Price (D) and Volume(1h) (E) columns
# for rows 2-51, including two row swaps (36<->37 and 44<->45).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.609.91"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "2.026.22"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.80"
$ws.Range("E5").Value = "  -2.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.634"
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.03"
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.388"
$ws.Range("E9").Value = "  +4.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.52"
$ws.Range("E10").Value = "  -2.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0787"
$ws.Range("E11").Value = "  +4.92%  "
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.900"
$ws.Range("E13").Value = "  -2.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.98"
$ws.Range("E14").Value = "  +16.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.29"
$ws.Range("E15").Value = "  -4.52%  "
$ws.Range("D16").Value = "2.326.48"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.50"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").Value = "2.030.72"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").Value = "36.601.60"
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.14"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "0.0₃0881"
$ws.Range("E21").Value = "  +2.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.36"
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.75"
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("E25").Value = "  -6.90%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.83"
$ws.Range("E27").Value = "  +3.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.137"
$ws.Range("E28").Value = "  +26.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "159.40"
$ws.Range("E29").Value = "  -3.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.19"
$ws.Range("E30").Value = "  +2.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.120"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.04"
$ws.Range("E32").Value = "  -1.59%  "
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0617"
$ws.Range("E34").Value = "  +1.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.48"
$ws.Range("E35").Value = "  -0.49%  "
$ws.Range("B36").Value = "THORChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.37"
$ws.Range("E36").Value = "  +8.91%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.36"
$ws.Range("E37").Value = "  -5.72%  "
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.16"
$ws.Range("E40").Value = "  +25.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0993"
$ws.Range("E41").Value = "  -4.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.25"
$ws.Range("E42").Value = "  +2.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.95"
$ws.Range("E43").Value = "  +1.99%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.13"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.93"
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0214"
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "93.49"
$ws.Range("E47").Value = "  -0.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.66"
$ws.Range("E48").Value = "  -2.21%  "
$ws.Range("D49").Value = "1.369.75"
$ws.Range("E49").Value = "  -3.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.90"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").Value = "2.214.77"
$ws.Range("E51").Value = "  +0.61%  "
